$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 5810.8667
$ws.Cells.Item(40, 9).Value = 4464.778
$ws.Cells.Item(40, 10).Value = 7830
$ws.Cells.Item(40, 11).Value = 4464.778
$ws.Cells.Item(40, 12).Value = 7830
$ws.Cells.Item(40, 13).Value = -4289.778
$ws.Cells.Item(40, 14).Value = -8180
$ws.Cells.Item(42, 8).Value = 171.08333
$ws.Cells.Item(42, 9).Value = 67.57143000000001
$ws.Cells.Item(42, 10).Value = 316
$ws.Cells.Item(42, 11).Value = 202.71429
$ws.Cells.Item(42, 12).Value = 948
$ws.Cells.Item(42, 13).Value = 27.28570999999999
$ws.Cells.Item(42, 14).Value = -1408
$ws.Cells.Item(94, 8).Value = 10692.733
$ws.Cells.Item(94, 9).Value = 10692.733
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 11).Value = 10692.733
$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 13).Value = -10241.733
$ws.Cells.Item(94, 14).ClearContents()
$ws.Cells.Item(121, 8).Value = 816.75
$ws.Cells.Item(121, 10).Value = 816.75
$ws.Cells.Item(121, 12).Value = 2450.25
$ws.Cells.Item(121, 14).Value = -5944.25
$ws.Cells.Item(129, 8).Value = 2580
$ws.Cells.Item(129, 10).Value = 2807
$ws.Cells.Item(129, 12).Value = 8421
$ws.Cells.Item(129, 14).Value = -18421
$ws.Cells.Item(131, 8).Value = 1320.1
$ws.Cells.Item(131, 9).Value = 986.2857
$ws.Cells.Item(131, 10).Value = 2099
$ws.Cells.Item(131, 11).Value = 2958.8571
$ws.Cells.Item(131, 12).Value = 6297
$ws.Cells.Item(131, 13).Value = 2081.1429
$ws.Cells.Item(131, 14).Value = -16377
$ws.Cells.Item(137, 8).Value = 2249.25
$ws.Cells.Item(137, 9).Value = 2249.25
$ws.Cells.Item(137, 11).Value = 6747.75
$ws.Cells.Item(137, 13).Value = -4197.75
$ws.Cells.Item(138, 8).Value = 4747.5
$ws.Cells.Item(138, 10).Value = 4946.6665
$ws.Cells.Item(138, 12).Value = 14839.9995
$ws.Cells.Item(138, 14).Value = -25119.9995

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(114, 8).Value = 63570.57
$ws.Cells.Item(114, 10).Value = 63570.57
$ws.Cells.Item(114, 12).Value = 63570.57
$ws.Cells.Item(114, 14).Value = -72248.57000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 856
$ws.Cells.Item(134, 9).Value = 856
$ws.Cells.Item(134, 11).Value = 2568
$ws.Cells.Item(134, 13).Value = -33

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 5485.5713
$ws.Cells.Item(22, 10).Value = 6250
$ws.Cells.Item(22, 12).Value = 6250
$ws.Cells.Item(22, 14).Value = -6950
$ws.Cells.Item(31, 8).Value = 7882
$ws.Cells.Item(31, 9).Value = 4857
$ws.Cells.Item(31, 10).Value = 9999.5
$ws.Cells.Item(31, 11).Value = 4857
$ws.Cells.Item(31, 12).Value = 9999.5
$ws.Cells.Item(31, 13).Value = -4562
$ws.Cells.Item(31, 14).Value = -10589.5
$ws.Cells.Item(32, 8).Value = 2010
$ws.Cells.Item(32, 9).Value = 2010
$ws.Cells.Item(32, 11).Value = 2010
$ws.Cells.Item(32, 13).Value = -1694
$ws.Cells.Item(34, 8).Value = 7882
$ws.Cells.Item(34, 9).Value = 4857
$ws.Cells.Item(34, 10).Value = 9999.5
$ws.Cells.Item(34, 11).Value = 4857
$ws.Cells.Item(34, 12).Value = 9999.5
$ws.Cells.Item(34, 13).Value = -4655
$ws.Cells.Item(34, 14).Value = -10403.5
$ws.Cells.Item(58, 8).Value = 3527.5
$ws.Cells.Item(58, 10).Value = 5555
$ws.Cells.Item(58, 12).Value = 5555
$ws.Cells.Item(58, 14).Value = -5961
$ws.Cells.Item(59, 8).Value = 70001.75
$ws.Cells.Item(59, 10).Value = 70001.75
$ws.Cells.Item(59, 12).Value = 70001.75
$ws.Cells.Item(59, 14).Value = -72291.75
$ws.Cells.Item(62, 8).Value = 7112.25
$ws.Cells.Item(62, 10).Value = 8483
$ws.Cells.Item(62, 12).Value = 8483
$ws.Cells.Item(62, 14).Value = -9731
$ws.Cells.Item(65, 8).Value = 7112.25
$ws.Cells.Item(65, 10).Value = 8483
$ws.Cells.Item(65, 12).Value = 42415
$ws.Cells.Item(65, 14).Value = -48655
$ws.Cells.Item(68, 8).Value = 243584
$ws.Cells.Item(68, 10).Value = 243584
$ws.Cells.Item(68, 12).Value = 243584
$ws.Cells.Item(68, 14).Value = -245082
$ws.Cells.Item(71, 8).Value = 243584
$ws.Cells.Item(71, 10).Value = 243584
$ws.Cells.Item(71, 12).Value = 730752
$ws.Cells.Item(71, 14).Value = -738240
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 14).ClearContents()
$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 14).ClearContents()
$ws.Cells.Item(125, 8).Value = 85000
$ws.Cells.Item(125, 10).Value = 85000
$ws.Cells.Item(125, 12).Value = 85000
$ws.Cells.Item(125, 14).Value = -89920
$ws.Cells.Item(132, 8).Value = 3813.9412
$ws.Cells.Item(132, 9).Value = 3295.5386
$ws.Cells.Item(132, 11).Value = 9886.6158
$ws.Cells.Item(132, 13).Value = -7356.6158
$ws.Cells.Item(136, 8).Value = 3527.5
$ws.Cells.Item(136, 10).Value = 5555
$ws.Cells.Item(136, 12).Value = 16665
$ws.Cells.Item(136, 14).Value = -21765
$ws.Cells.Item(141, 8).Value = 61633.332
$ws.Cells.Item(141, 10).Value = 77450
$ws.Cells.Item(141, 12).Value = 77450
$ws.Cells.Item(141, 14).Value = -87810

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(60, 8).Value = 620.9048
$ws.Cells.Item(60, 9).Value = 279.66666
$ws.Cells.Item(60, 11).Value = 838.9999799999999
$ws.Cells.Item(60, 13).Value = -587.9999799999999
$ws.Cells.Item(109, 8).Value = 275606.5
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 1943.6364
$ws.Cells.Item(113, 10).Value = 1999
$ws.Cells.Item(113, 12).Value = 5997
$ws.Cells.Item(113, 14).Value = -10337

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 6156.0713
$ws.Cells.Item(113, 9).Value = 3455
$ws.Cells.Item(113, 11).Value = 3455
$ws.Cells.Item(113, 13).Value = -1285
$ws.Cells.Item(132, 8).Value = 97024.336
$ws.Cells.Item(132, 9).Value = 138461.5
$ws.Cells.Item(132, 11).Value = 415384.5
$ws.Cells.Item(132, 13).Value = -412854.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 699.0769
$ws.Cells.Item(22, 9).Value = 697.25
$ws.Cells.Item(22, 10).Value = 699.8889
$ws.Cells.Item(22, 11).Value = 697.25
$ws.Cells.Item(22, 12).Value = 699.8889
$ws.Cells.Item(22, 13).Value = -402.25
$ws.Cells.Item(22, 14).Value = -1289.8889
$ws.Cells.Item(27, 8).Value = 699.0769
$ws.Cells.Item(27, 9).Value = 697.25
$ws.Cells.Item(27, 10).Value = 699.8889
$ws.Cells.Item(27, 11).Value = 697.25
$ws.Cells.Item(27, 12).Value = 699.8889
$ws.Cells.Item(27, 13).Value = -590.25
$ws.Cells.Item(27, 14).Value = -913.8889
$ws.Cells.Item(46, 8).Value = 6252.7144
$ws.Cells.Item(46, 9).Value = 1909.5
$ws.Cells.Item(46, 10).Value = 7990
$ws.Cells.Item(46, 11).Value = 1909.5
$ws.Cells.Item(46, 12).Value = 7990
$ws.Cells.Item(46, 13).Value = -1721.5
$ws.Cells.Item(46, 14).Value = -8366
$ws.Cells.Item(68, 8).Value = 6728.5713
$ws.Cells.Item(68, 10).Value = 8700
$ws.Cells.Item(68, 12).Value = 8700
$ws.Cells.Item(68, 14).Value = -10198
$ws.Cells.Item(71, 8).Value = 6728.5713
$ws.Cells.Item(71, 10).Value = 8700
$ws.Cells.Item(71, 12).Value = 43500
$ws.Cells.Item(71, 14).Value = -50988
$ws.Cells.Item(82, 8).Value = 3615.5881
$ws.Cells.Item(82, 9).Value = 2618.111
$ws.Cells.Item(82, 10).Value = 4737.75
$ws.Cells.Item(82, 11).Value = 2618.111
$ws.Cells.Item(82, 12).Value = 4737.75
$ws.Cells.Item(82, 13).Value = -2257.111
$ws.Cells.Item(82, 14).Value = -5459.75
$ws.Cells.Item(85, 8).Value = 3615.5881
$ws.Cells.Item(85, 9).Value = 2618.111
$ws.Cells.Item(85, 10).Value = 4737.75
$ws.Cells.Item(85, 11).Value = 2618.111
$ws.Cells.Item(85, 12).Value = 4737.75
$ws.Cells.Item(85, 13).Value = -1370.111
$ws.Cells.Item(85, 14).Value = -7233.75
$ws.Cells.Item(125, 8).Value = 65000
$ws.Cells.Item(125, 10).Value = 65000
$ws.Cells.Item(125, 12).Value = 65000
$ws.Cells.Item(125, 14).Value = -74840

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 14).ClearContents()
$ws.Cells.Item(109, 8).Value = 99999
$ws.Cells.Item(109, 10).Value = 99999
$ws.Cells.Item(109, 12).Value = 99999
$ws.Cells.Item(109, 14).Value = -102773
